# ----------------------------------------------------------------------------
# Applies the "add few more service to existing contractMaker" edit:
#
#  1. In the Party A (甲方) paragraph, the run of three spaces that sits right
#     before "${order_company_name}${order_company_tax_ref}" is shortened to a
#     single space, and the run of six trailing spaces at the very end of that
#     paragraph (after the other copy of the three-space run) is removed.
#  2. In the Party B (乙方) paragraph, one of the two adjacent two-space runs
#     that sit right before "${firm_name}" is removed (collapsing four spaces
#     down to two).
#  3. The "_GoBack" bookmark is moved from the end of clause 2 ("...传真件与
#     合同原件具有同等效力。") to sit right after the "firm_name" text (and
#     right before the closing "}") in the Party B paragraph.
#
# NOTE: text containing "${...}" must be single-quoted in PowerShell --
# double quotes would trigger variable-expansion syntax.
# ----------------------------------------------------------------------------

$d = $word.ActiveDocument

# --- Party A paragraph -------------------------------------------------
# Locate the merged template marker run.
$markerA = $d.Content
$okA = $markerA.Find.Execute('${order_company_name}${order_company_tax_ref}', `
        $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if (-not $okA) {
    throw "Could not find the order_company_name/order_company_tax_ref marker"
}

# The run directly before the marker is the three-space run -> shrink to one space.
$preSpaces = $d.Range($markerA.Start - 3, $markerA.Start)
$preSpaces.Text = " "

# After shrinking, the marker's end shifted back by 2 characters.
$afterMarkerStart = $markerA.End - 2

# Following the marker there are two runs: "   " (3 spaces, kept) then
# "      " (6 spaces, removed). Drop the trailing six spaces.
$trailingSpaces = $d.Range($afterMarkerStart + 3, $afterMarkerStart + 9)
$trailingSpaces.Text = ""

# --- Party B paragraph ---------------------------------------------------
# Locate "firm_name" (the bare text, not the surrounding ${ }).
$firmName = $d.Content
$okB = $firmName.Find.Execute('firm_name', $true, $false, $false, $false, `
        $false, $true, 1, $false, "", 0)

if (-not $okB) {
    throw "Could not find the firm_name marker"
}

# Immediately before "${firm_name}" there are two adjacent two-space runs
# (four spaces total). Remove the second (closer) one, leaving two spaces.
$innerSpaces = $d.Range($firmName.Start - 4, $firmName.Start - 2)
$innerSpaces.Text = ""

# Re-find "firm_name" now that the deletion above shifted offsets.
$firmName2 = $d.Content
$null = $firmName2.Find.Execute('firm_name', $true, $false, $false, $false, `
        $false, $true, 1, $false, "", 0)

# --- Move the "_GoBack" bookmark -----------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

$newBookmarkRange = $d.Range($firmName2.End, $firmName2.End)
$d.Bookmarks.Add("_GoBack", $newBookmarkRange)

Write-Host "Edit applied successfully"
